$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy existing row formatting down onto the new rows (241-260) ---
# Rows with only A (number, style 1) + B (text, style 2) -> mimic row 240
$ws.Range("A240:B240").Copy() | Out-Null
$ws.Range("A241:B246").PasteSpecial(-4122) | Out-Null
$ws.Range("A248:B252").PasteSpecial(-4122) | Out-Null
$ws.Range("A254:B257").PasteSpecial(-4122) | Out-Null
$ws.Range("A259:B260").PasteSpecial(-4122) | Out-Null

# Rows with A (number, style 1) + B (text, style 2) + C (text, style 2) -> mimic row 43
$ws.Range("A43:C43").Copy() | Out-Null
$ws.Range("A247:C247").PasteSpecial(-4122) | Out-Null
$ws.Range("A253:C253").PasteSpecial(-4122) | Out-Null

# Row 258: A (number, style1) + B (text, style2) + C (text, style1) + D (text, style2)
$ws.Range("A43:B43").Copy() | Out-Null
$ws.Range("A258:B258").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Copy() | Out-Null
$ws.Range("C258").PasteSpecial(-4122) | Out-Null
$ws.Range("D43").Copy() | Out-Null
$ws.Range("D258").PasteSpecial(-4122) | Out-Null

# --- Fill in values ---
$ws.Range("A241").Value = 2
$ws.Range("B241").Value = "飞鸟"

$ws.Range("A242").Value = 2
$ws.Range("B242").Value = "Caesar"

$ws.Range("A243").Value = 2
$ws.Range("B243").Value = "张凯"

$ws.Range("A244").Value = 2
$ws.Range("B244").Value = "宅男阿海"

$ws.Range("A245").Value = 2
$ws.Range("B245").Value = "程蝶衣"

$ws.Range("A246").Value = 2
$ws.Range("B246").Value = "八千"

$ws.Range("A247").Value = 1
$ws.Range("B247").Value = "b 安贝慧"
$ws.Range("C247").Value = "  "

$ws.Range("A248").Value = 2
$ws.Range("B248").Value = "Alex"

$ws.Range("A249").Value = 2
$ws.Range("B249").Value = "呼噜呼噜"

$ws.Range("A250").Value = 2
$ws.Range("B250").Value = "解无明"

$ws.Range("A251").Value = 2
$ws.Range("B251").Value = "千鹤"

$ws.Range("A252").Value = 2
$ws.Range("B252").Value = "zeqing"

$ws.Range("A253").Value = 1
$ws.Range("B253").Value = "Si_X:"
$ws.Range("C253").Value = "社畜一枚，玩着你们的游戏睡着了，我做了一个梦，这个梦真的很美，谢谢"

$ws.Range("A254").Value = 2
$ws.Range("B254").Value = "不喜欢甜的"

$ws.Range("A255").Value = 2
$ws.Range("B255").Value = "晓危"

$ws.Range("A256").Value = 2
$ws.Range("B256").Value = "毛虫哥哥"

$ws.Range("A257").Value = 2
$ws.Range("B257").Value = "佚名"

$ws.Range("A258").Value = 1
$ws.Range("B258").Value = "听风忆雪"
$ws.Range("C258").Value = "搞死我了"
$ws.Range("D258").Value = "#4169E1"

$ws.Range("A259").Value = 2
$ws.Range("B259").Value = "老郑"

$ws.Range("A260").Value = 2
$ws.Range("B260").Value = "太帅很苦恼"

# --- Update the saved view state to roughly match the author's last position ---
$ws.Range("I233").Select() | Out-Null
